$wb = $excel.ActiveWorkbook

# 1. Rename the "MFO fleet" sheet to "MDO fleet".
#    Excel automatically rewrites every formula reference to the renamed
#    sheet (e.g. in the "total" sheet), so no further formula edits needed.
$mfo = $wb.Worksheets.Item("MFO fleet")
$mfo.Name = "MDO fleet"

# 2. Update the fuel-code labels on the "ships" sheet from *_MFO to *_MDO.
$ships = $wb.Worksheets.Item("ships")
$ships.Range("B3").Value = """T_MDO"""
$ships.Range("B4").Value = """B_MDO"""
$ships.Range("B5").Value = """G_MDO"""
$ships.Range("B6").Value = """C_MDO"""
$ships.Range("B7").Value = """O_MDO"""

# 3. Move the active tab / selection from "MDO fleet" to "ships", with B8
#    selected as the active cell there.
$ships.Activate()
$ships.Range("B8").Select()
